$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume(1h) (E) columns to remain text so numeric-looking values are not auto-converted
$ws.Range("D2:E51").NumberFormat = "@"

# Row 39/40: VeChain and Aptos swap positions (rank index in column A stays fixed)
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01989"
$ws.Range("E39").Value = "  -3.01%  "

$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "10.17"
$ws.Range("E40").Value = "  -0.56%  "

# Remaining per-row Price/Volume(1h) updates
$ws.Range("D2").Value = "20.265.77"
$ws.Range("E2").Value = "  +2.28%  "
$ws.Range("D3").Value = "1.442.24"
$ws.Range("E3").Value = "  +3.22%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "0.9146"
$ws.Range("E5").Value = "  -8.85%  "
$ws.Range("D6").Value = "274.79"
$ws.Range("E6").Value = "  +0.73%  "
$ws.Range("D7").Value = "0.3630"
$ws.Range("E7").Value = "  -0.95%  "
$ws.Range("D8").Value = "0.3080"
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("D9").Value = "39.24"
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("E10").Value = "  +2.36%  "
$ws.Range("D11").Value = "0.06512"
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").Value = "0.9984"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").Value = "5.353"
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("D14").Value = "17.53"
$ws.Range("E14").Value = "  +2.26%  "
$ws.Range("D15").Value = "6.043"
$ws.Range("E15").Value = "  -1.12%  "
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Value = "1.439.82"
$ws.Range("E17").Value = "  +2.55%  "
$ws.Range("D18").Value = "0.9299"
$ws.Range("E18").Value = "  -7.19%  "
$ws.Range("D19").Value = "0.05618"
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("D20").Value = "67.60"
$ws.Range("E20").Value = "  -3.44%  "
$ws.Range("D21").Value = "5.415"
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("D22").Value = "14.23"
$ws.Range("E22").Value = "  -2.53%  "
$ws.Range("D23").Value = "10.82"
$ws.Range("E23").Value = "  -1.08%  "
$ws.Range("D24").Value = "2.230"
$ws.Range("E24").Value = "  -1.77%  "
$ws.Range("D25").Value = "20.287.11"
$ws.Range("E25").Value = "  +2.11%  "
$ws.Range("D26").Value = "137.67"
$ws.Range("E26").Value = "  +2.07%  "
$ws.Range("D27").Value = "2.078"
$ws.Range("E27").Value = "  -5.63%  "
$ws.Range("D28").Value = "16.93"
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("D29").Value = "1.590.94"
$ws.Range("E29").Value = "  +1.83%  "
$ws.Range("D30").Value = "110.22"
$ws.Range("E30").Value = "  +1.28%  "
$ws.Range("D31").Value = "3.958"
$ws.Range("E31").Value = "  -2.83%  "
$ws.Range("D32").Value = "0.8016"
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("D33").Value = "4.850"
$ws.Range("E33").Value = "  -7.45%  "
$ws.Range("D34").Value = "0.07665"
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("D35").Value = "1.462"
$ws.Range("E35").Value = "  +1.36%  "
$ws.Range("D36").Value = "0.05801"
$ws.Range("E36").Value = "  +1.19%  "
$ws.Range("D37").Value = "4.671"
$ws.Range("E37").Value = "  -2.45%  "
$ws.Range("D38").Value = "1.132"
$ws.Range("E38").Value = "  +3.81%  "
$ws.Range("D41").Value = "0.1852"
$ws.Range("E41").Value = "  -2.03%  "
$ws.Range("D42").Value = "0.9263"
$ws.Range("E42").Value = "  -7.45%  "
$ws.Range("D43").Value = "7.074"
$ws.Range("E43").Value = "  -15.18%  "
$ws.Range("D44").Value = "0.5207"
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("D45").Value = "3.482"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("D46").Value = "11.87"
$ws.Range("E46").Value = "  -3.37%  "
$ws.Range("D47").Value = "116.78"
$ws.Range("E47").Value = "  +5.16%  "
$ws.Range("D48").Value = "0.5106"
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("D49").Value = "1.729"
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("D50").Value = "0.06420"
$ws.Range("E50").Value = "  +4.69%  "
$ws.Range("D51").Value = "0.9758"
$ws.Range("E51").Value = "  -2.53%  "